$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target contents for rows 2-10 (columns A-T), following Dr Hou's advice:
# a new "ECs" cluster is added to the sending/target cluster categories and
# every Sending/Target-cluster combination (ECs, FAPs, sCs) is recomputed.
$data = @(
    @("ECs","Fat4","Dchs1","ECs",2,0.6666666666666666,18.40105966666667,55.203179,0.6177904080289686,0.6177904080289685,2,0.6666666666666666,7.519968333333334,22.559905,0.483567271060424,0.483567271060424,138.3753859931105,1245.378473937995,0.2987432216978742,0.2987432216978741),
    @("ECs","Fat4","Dchs1","FAPs",2,0.6666666666666666,18.40105966666667,55.203179,0.6177904080289686,0.6177904080289685,3,1,7.006278666666667,21.018836,0.4505347502742851,0.4505347502742851,128.9229517866271,1160.306566079644,0.2783360472031801,0.27833604720318),
    @("ECs","Fat4","Dchs1","sCs",2,0.6666666666666666,18.40105966666667,55.203179,0.6177904080289686,0.6177904080289685,3,1,1.024781333333333,3.074344,0.06589797866529083,0.06589797866529083,18.85706245995289,169.713562139576,0.04071113912791429,0.04071113912791428),
    @("FAPs","Fat4","Dchs1","ECs",3,1,8.741643666666667,26.224931,0.2934887286658174,0.2934887286658173,2,0.6666666666666666,7.519968333333334,22.559905,0.483567271060424,0.483567271060424,65.73688355461722,591.631951991555,0.1419215436079225,0.1419215436079225),
    @("FAPs","Fat4","Dchs1","FAPs",3,1,8.741643666666667,26.224931,0.2934887286658174,0.2934887286658173,3,1,7.006278666666667,21.018836,0.4505347502742851,0.4505347502742851,61.24639153336845,551.2175238003159,0.1322268710777715,0.1322268710777714),
    @("FAPs","Fat4","Dchs1","sCs",3,1,8.741643666666667,26.224931,0.2934887286658174,0.2934887286658173,3,1,1.024781333333333,3.074344,0.06589797866529083,0.06589797866529083,8.958273252251555,80.624459270264,0.01934031398012336,0.01934031398012336),
    @("sCs","Fat4","Dchs1","ECs",3,1,2.642575666666667,7.927727,0.0887208633052142,0.08872086330521417,2,0.6666666666666666,7.519968333333334,22.559905,0.483567271060424,0.483567271060424,19.87208533177056,178.848767985935,0.04290250575462733,0.04290250575462732),
    @("sCs","Fat4","Dchs1","FAPs",3,1,2.642575666666667,7.927727,0.0887208633052142,0.08872086330521417,3,1,7.006278666666667,21.018836,0.4505347502742851,0.4505347502742851,18.51462151841911,166.631593665772,0.03997183199333366,0.03997183199333365),
    @("sCs","Fat4","Dchs1","sCs",3,1,2.642575666666667,7.927727,0.0887208633052142,0.08872086330521417,3,1,1.024781333333333,3.074344,0.06589797866529083,0.06589797866529083,2.708062215120889,24.372559936088,0.00584652555725319,0.005846525557253188)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
